# Rework the counters_summary sheet: the old per-dimension metric columns
# (UNIQUENESS / UNIQUENESSCOMPOSITE / FORMATCONSISTENCY / ... / CONFORMANCE)
# are replaced by a generated "metric + metric SCORE" column pair for every
# data-quality dimension (completeness-mandatory, completeness-optional,
# precision, business-rule compliance, metadata compliance, uniqueness,
# non-redundancy, semantic consistency, value consistency, format
# consistency). Each SCORE column is the raw counter divided by 24.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$allRows = @(
    @('attribute', 'COMPLETENESSMANDATORY', 'COMPLETENESSMANDATORY SCORE', 'COMPLETENESSOPTIONAL', 'COMPLETENESSOPTIONAL SCORE', 'PRECISION', 'PRECISION SCORE', 'BUSINESSRULECOMPLIANCE', 'BUSINESSRULECOMPLIANCE SCORE', 'METADATACOMPLIANCE', 'METADATACOMPLIANCE SCORE', 'UNIQUENESS', 'UNIQUENESS SCORE', 'NONREDUNDANCY', 'NONREDUNDANCY SCORE', 'SEMANTICCONSISTENCY', 'SEMANTICCONSISTENCY SCORE', 'VALUECONSISTENCY', 'VALUECONSISTENCY SCORE', 'FORMATCONSISTENCY', 'FORMATCONSISTENCY SCORE'),
    @('caseNumber', 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 1, 0.041667, 0, 0, 0, 0, 0, 0, 8, 0.333333),
    @('sentence', 0, 0, 0, 0, 0, 0, 0, 0, 3, 0.125, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0),
    @('registrationDate', 12, 0.5, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 20, 0.833333),
    @('category', 0, 0, 0, 0, 0, 0, 0, 0, 7, 0.291667, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0),
    @('firstName', 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0),
    @('middleName', 5, 0.208333, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0),
    @('lastName', 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0),
    @('age', 0, 0, 0, 0, 0, 0, 0, 0, 1, 0.041667, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0),
    @('sex', 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0),
    @('race', 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0),
    @('district', 4, 0.166667, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0),
    @('post', 4, 0.166667, 0, 0, 0, 0, 0, 0, 4, 0.166667, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0),
    @('neighborhood', 4, 0.166667, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0),
    @('Location 1', 4, 0.166667, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
)

for ($r = 0; $r -lt $allRows.Length; $r++) {
    $rowValues = $allRows[$r]
    for ($c = 0; $c -lt $rowValues.Length; $c++) {
        $ws.Cells.Item($r + 1, $c + 1).Value = $rowValues[$c]
    }
}

Write-Host "Rewrote counters_summary header/data grid (A1:U15) with the new per-dimension score columns."
